$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01134666666666667
$ws.Range("H2").Value = 0.03404
$ws.Range("I2").Value = 0.001209510404472147
$ws.Range("J2").Value = 0.001209510404472147
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 0.005904211524444446
$ws.Range("R2").Value = 0.05313790372
$ws.Range("S2").Value = 0.000004966170348174363
$ws.Range("T2").Value = 0.000004966170348174365

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01134666666666667
$ws.Range("H3").Value = 0.03404
$ws.Range("I3").Value = 0.001209510404472147
$ws.Range("J3").Value = 0.001209510404472147
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 1.202329408128889
$ws.Range("R3").Value = 10.82096467316
$ws.Range("S3").Value = 0.00101130737451849
$ws.Range("T3").Value = 0.00101130737451849

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01134666666666667
$ws.Range("H4").Value = 0.03404
$ws.Range("I4").Value = 0.001209510404472147
$ws.Range("J4").Value = 0.001209510404472147
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 0.2297366408
$ws.Range("R4").Value = 2.0676297672
$ws.Range("S4").Value = 0.0001932368596054828
$ws.Range("T4").Value = 0.0001932368596054828

# Row 5
$ws.Range("I5").Value = 0.8865539289740954
$ws.Range("J5").Value = 0.8865539289740952
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 4.327703097994224
$ws.Range("R5").Value = 38.94932788194801
$ws.Range("S5").Value = 0.003640132253388997
$ws.Range("T5").Value = 0.003640132253388998

# Row 6
$ws.Range("I6").Value = 0.8865539289740954
$ws.Range("J6").Value = 0.8865539289740952
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("S6").Value = 0.7412739261810052
$ws.Range("T6").Value = 0.7412739261810052

# Row 7
$ws.Range("I7").Value = 0.8865539289740954
$ws.Range("J7").Value = 0.8865539289740952
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 168.39369118072
$ws.Range("R7").Value = 1515.54322062648
$ws.Range("S7").Value = 0.1416398705397011
$ws.Range("T7").Value = 0.1416398705397011

# Row 8
$ws.Range("G8").Value = 1.052914333333334
$ws.Range("H8").Value = 3.158743
$ws.Range("I8").Value = 0.1122365606214325
$ws.Range("J8").Value = 0.1122365606214325
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5203476666666668
$ws.Range("N8").Value = 1.561043
$ws.Range("O8").Value = 0.004105934376266647
$ws.Range("P8").Value = 0.004105934376266647
$ws.Range("Q8").Value = 0.547881516549889
$ws.Range("R8").Value = 4.930933648949001
$ws.Range("S8").Value = 0.0004608359525294752
$ws.Range("T8").Value = 0.0004608359525294753

# Row 9
$ws.Range("G9").Value = 1.052914333333334
$ws.Range("H9").Value = 3.158743
$ws.Range("I9").Value = 0.1122365606214325
$ws.Range("J9").Value = 0.1122365606214325
$ws.Range("O9").Value = 0.8361295370252257
$ws.Range("P9").Value = 0.8361295370252259
$ws.Range("Q9").Value = 111.5701998126108
$ws.Range("R9").Value = 1004.131798313497
$ws.Range("S9").Value = 0.09384430346970206
$ws.Range("T9").Value = 0.09384430346970206

# Row 10
$ws.Range("G10").Value = 1.052914333333334
$ws.Range("H10").Value = 3.158743
$ws.Range("I10").Value = 0.1122365606214325
$ws.Range("J10").Value = 0.1122365606214325
$ws.Range("M10").Value = 20.24706
$ws.Range("N10").Value = 60.74118
$ws.Range("O10").Value = 0.1597645285985076
$ws.Range("P10").Value = 0.1597645285985076
$ws.Range("Q10").Value = 21.31841968186
$ws.Range("R10").Value = 191.86577713674
$ws.Range("S10").Value = 0.01793142119920099
$ws.Range("T10").Value = 0.01793142119920099
